# Update handback status timestamps ("Generate Report for Handback")
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for d75e65e5-...md (row 3, col G)
$wsOverview.Range("G3").Value = "2016-08-18 02:39:41"

# zh-cn sheet: d75e65e5-... row (row 3)
#   Correspond Handoff Datetime (col H)
$wsZhCn.Range("H3").Value = "2016-08-18 02:39:36"
#   Correspond Handback DateTime (col K)
$wsZhCn.Range("K3").Value = "2016-08-18 02:39:52"

# de-de sheet: d75e65e5-... row (row 3)
#   Correspond Handback DateTime (col K)
$wsDeDe.Range("K3").Value = "2016-08-18 02:40:02"
